$d = $word.ActiveDocument

# The document contains two <id>...</id> tag sequences, each split across
# three separate runs: "<id>", the literal id value, and "</id>". Merge
# each trio back into a single run (matching the formatting already
# carried by the opening "<id>" run) by finding the full tag text as it
# reads across the runs and replacing it with itself.

$d.Content.Find.Execute("<id>p106r_1</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p106r_1</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p106r_2</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p106r_2</id>", 2) | Out-Null
